$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Participants")

### ---- Set raw values first (while D48/D49 are still plain/General format,
###      so the numeric Code values are not coerced to text by the table's
###      Text-formatted "Code" column) ----

$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()

# Names/IDs for the three new rows are entered first (matches the original
# authoring sequence, and therefore the shared-string insertion order).
$ws.Range("B48").Value = "Naroa Borrajeros Elósegui"
$ws.Range("B49").Value = "Thais Garcia Sevilla"
$ws.Range("B50").Value = "Unai Roca"

# Row 48: new participant 47 - Naroa Borrajeros Elosegui
$ws.Range("A48").Value = 47
$ws.Range("C48").Value = 6020
$ws.Range("D48").Value = 470104
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 4
$ws.Range("G48").Value = "C04"
$ws.Range("H48").Value = 45317
$ws.Range("I48").Value = "A bit sleepy. Can do the caja."
$ws.Range("J48").Value = 36284
$ws.Range("K48").Value = "F"

# Row 49: new participant 48 - Thais Garcia Sevilla
$ws.Range("A49").Value = 48
$ws.Range("C49").Value = 10038
$ws.Range("D49").Value = 480204
$ws.Range("E49").Value = 2
$ws.Range("F49").Value = 4
$ws.Range("G49").Value = "C01"
$ws.Range("H49").Value = 45317
$ws.Range("I49").Value = "Has a cut on her right index. No issue in typing. Took off the bandage. Can do the caja."
$ws.Range("J49").Value = 33386
$ws.Range("K49").Value = "F"

# Row 50: participant 49 also gets a Name/ID - Unai Roca
$ws.Range("C50").Value = 6477

# Row 51: new trailing spacer row.
$ws.Range("I51").Value = " "

### ---- Now copy formatting from matching template cells (values untouched) ----

# Row 20 has the exact style pattern needed for columns A,B,C,D,E,F,G,H,J,K of rows 48 & 49
# (I is intentionally left with the default/no style, same as the template).
$ws.Range("A20:K20").Copy() | Out-Null
$ws.Range("A48:K48").PasteSpecial(-4122) | Out-Null
$ws.Range("A49:K49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 50 B/C take the bold "s=10" style (same style already used for column I elsewhere, e.g. I18).
$ws.Range("I18").Copy() | Out-Null
$ws.Range("B50:C50").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 51: D51/G51 reuse the empty "s=1" style that D50/G50 already carry.
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4122) | Out-Null
$ws.Range("G50").Copy() | Out-Null
$ws.Range("G51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# K51 inherits K50's current "s=22" placeholder style, then K50 itself is cleared entirely.
$ws.Range("K50").Copy() | Out-Null
$ws.Range("K51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("K50").Clear()

# Resize the table to include the newly added row 51.
$ws.ListObjects.Item("Table1").Resize($ws.Range("A1:O51"))

$ws.Range("A41").Select()
$ws.Range("I50").Select()
